# Weekly fruit/vegetable price update:
# Insert a new daily price record as row 70 (pushing the existing
# rows 70-112 down to 71-113) on the "Cebollín" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 70:112 down one row to make room for the new record.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Cells.Item(70, 1).Value = 11
$ws.Cells.Item(70, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(70, 3).Value = "Bíobío"
$ws.Cells.Item(70, 4).Value = 45072
$ws.Cells.Item(70, 5).Value = 8
$ws.Cells.Item(70, 6).Value = 100112037
$ws.Cells.Item(70, 7).Value = "Cebollín"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 110
$ws.Cells.Item(70, 11).Value = 5000
$ws.Cells.Item(70, 12).Value = 5500
$ws.Cells.Item(70, 13).Value = 5273
$ws.Cells.Item(70, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(70, 16).Value = 146
$ws.Cells.Item(70, 17).Value = 36
$ws.Cells.Item(70, 18).Value = "Hortaliza"
